$d = $word.ActiveDocument
$rng = $d.Content
$rng.Collapse(0)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Technical Side of CMS: Using PHP CodeIgniter as a framework:</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Creating a brand-new project:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">1.) on config, update </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>baseurl</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> to localhost/&lt;</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>foldername</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>&gt;/</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">2.) on routes, update </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>default_controller</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> to landing controller &lt;</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>nameOfController</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>&gt;</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>3.) Create a controller for the default, or copy paste welcome page.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>4.) Create the view of that default controller.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">5.) Create the database in phpMyAdmin, update </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>database.php</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, input database information</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>6.) In model (or created model), update database name</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>7.) rinse and repeat for different processes…</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>How to search for templates:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>- google ERP System Template HTML</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>- google Bootstrap admin templates</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Bugs (Backlog):</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">- </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>img</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> user logo is not showing when in </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>viewMember</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> or </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>viewMinistry</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Church Management System – Project Deadline and Milestones:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Ministry Section:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">What needs to be done? </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">- Add </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">a member </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>to ministry modal</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve"> on members</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>’</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve"> page</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve"> (members and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>viewMember</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve"> page for the select box is not working properly)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve"> and also save button.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">- Add member button to link to </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>members’</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve"> page</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">- Print </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">the members of the ministry in a table </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">to </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">ministry </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>view</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>- Edit ministry information</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>- Ministry description</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>, a select statement when Edit Ministry is clicked.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>Review Points (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>Aljon</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>):</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>Ministry Page:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="70AD47" w:themeColor="accent6"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="70AD47" w:themeColor="accent6"/>
        </w:rPr>
        <w:t>Delete Ministry</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="70AD47" w:themeColor="accent6"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="70AD47" w:themeColor="accent6"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="70AD47" w:themeColor="accent6"/>
        </w:rPr>
        <w:t>- After deleting a ministry, it is not displaying the members after redirecting.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="70AD47" w:themeColor="accent6"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="70AD47" w:themeColor="accent6"/>
        </w:rPr>
        <w:t>Add Ministry</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="70AD47" w:themeColor="accent6"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="70AD47" w:themeColor="accent6"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="70AD47" w:themeColor="accent6"/>
        </w:rPr>
        <w:tab/>
        <w:t>- After adding too many ministries, the scrollbar is not working.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>Add Members Tab on the Left Pane</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>When editing a Ministry, redirect to the specific Ministry edited</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>Change function names in the model to get/set/insert/delete for uniformity</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">Change function name in the controller to view/edit etc. instead of </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>viewMinistry</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve"> etc.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">Change the position of the Ministry Description imitating how the Network is displayed in the members </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>page.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>s</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>Remove colon on the ministry members</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>Retrieve Attendance Functionality on the Windows laptop</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Deadline:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t xml:space="preserve">June 8, 2022 (Wednesday) </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Events and Attendance Section:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>- Do delete event button to delete an event.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
        </w:rPr>
        <w:t>- Initialize Status to Open.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">- </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Adhoc</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> task when an attendee is not a member yet.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>- Create a “Non-member” button on the 2</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="superscript"/>
        </w:rPr>
        <w:t>nd</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> tab to input name (required). Then optional fields are network, address, and contact number.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> When clicked, it should pop out a modal </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">with these fields, then add it to the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ci_event_attendee</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, and print on the next tab.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">- </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>So</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> you have to create a new table called “</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ci_evt_non_member</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">” with fields (name, address, contact number, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>first_timer_flg</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> (yes or no values only)).</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>- Extract to pdf or excel.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">- Add another </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>tab(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>third tab) on the wizard to print attendees on the event</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> on a table</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Deadline</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>?</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>?</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Finance Section:</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Learnings along the way of developing this project:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>- Think in a business perspective. What do the admin users want and need to see for their purposes?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>- When debugging, find the cause of the problem first and start there.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>- Always comment so other programmers can understand the code and be able to help you fix the problem faster.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>- Know the design of the system first before setting the deadline.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>- The skills I need to acquire is learning how to find the problem, put them into tiny pieces, and fix it.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>- The starting point of debugging is to comment codes or print codes line by line, see if the value being passed is correct, where does the code explodes, etc.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Once you spot where the problem is, work towards fixing that.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>- In models, it should only contain getters and setters. Naming convention is important.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rng.InsertXML($xml)
Write-Output "Inserted new content."
